$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 30.26572128279039
$ws.Range("D2").Value = 0.4257212827903949
$ws.Range("E2").Value = 0.1812386106206994
$ws.Range("C3").Value = 30.13497992549404
$ws.Range("D3").Value = 0.3249799254940449
$ws.Range("E3").Value = 0.105611951974115
$ws.Range("C4").Value = 29.88282032601017
$ws.Range("D4").Value = -0.03717967398983291
$ws.Range("E4").Value = 0.001382328157990258
$ws.Range("C5").Value = 29.89411918165291
$ws.Range("D5").Value = -0.08588081834708561
$ws.Range("E5").Value = 0.007375514959965116
$ws.Range("C6").Value = 29.94241012696325
$ws.Range("D6").Value = -0.09758987303674971
$ws.Range("E6").Value = 0.009523783319328927
$ws.Range("C7").Value = 30.00820204159942
$ws.Range("D7").Value = -0.2017979584005793
$ws.Range("E7").Value = 0.04072241601464195
$ws.Range("C8").Value = 30.14504543960305
$ws.Range("D8").Value = -0.0749545603969537
$ws.Range("E8").Value = 0.005618186124300579
$ws.Range("C9").Value = 30.20597661290189
$ws.Range("D9").Value = -0.1740233870981065
$ws.Range("E9").Value = 0.03028413925709742
$ws.Range("C10").Value = 30.49339099702649
$ws.Range("D10").Value = 0.05339099702649008
$ws.Range("E10").Value = 0.002850598563482673
$ws.Range("C11").Value = 30.50756586496083
$ws.Range("D11").Value = 0.02756586496083102
$ws.Range("E11").Value = 0.0007598769110387715
$ws.Range("C12").Value = 30.51637331853345
$ws.Range("D12").Value = -0.1736266814665548
$ws.Range("E12").Value = 0.03014622451708848
$ws.Range("C13").Value = 30.73448260021426
$ws.Range("D13").Value = -0.01551739978573963
$ws.Range("E13").Value = 0.0002407896961104722
$ws.Range("C14").Value = 30.86336916496587
$ws.Range("D14").Value = -0.07663083503413404
$ws.Range("E14").Value = 0.005872284878028665
$ws.Range("C15").Value = 31.03914848366619
$ws.Range("D15").Value = 0.08914848366618955
$ws.Range("E15").Value = 0.007947452139980864
$ws.Range("C16").Value = 31.21170990337458
$ws.Range("D16").Value = 0.1917099033745835
$ws.Range("E16").Value = 0.03675268705189214
$ws.Range("C17").Value = 31.27962889131896
$ws.Range("D17").Value = 0.1596288913189596
$ws.Range("E17").Value = 0.02548138294372023
$ws.Range("C18").Value = 31.11440397279001
$ws.Range("D18").Value = -0.1655960272099897
$ws.Range("E18").Value = 0.02742204422773164
$ws.Range("C19").Value = 30.89072679564211
$ws.Range("D19").Value = -0.4892732043578896
$ws.Range("E19").Value = 0.2393882685026372
$ws.Range("C20").Value = 31.11169872854706
$ws.Range("D20").Value = -0.4683012714529333
$ws.Range("E20").Value = 0.219306080844434
$ws.Range("C21").Value = 31.65519723864403
$ws.Range("D21").Value = 0.005197238644029767
$ws.Range("E21").Value = 0.00002701128952299637
$ws.Range("C22").Value = 32.54992960015406
$ws.Range("D22").Value = 0.6699296001540596
$ws.Range("E22").Value = 0.4488056691625781
$ws.Range("C23").Value = 32.57320848691319
$ws.Range("D23").Value = 0.293208486913187
$ws.Range("E23").Value = 0.08597121679792057
$ws.Range("C24").Value = 32.77063300299165
$ws.Range("D24").Value = 0.3206330029916487
$ws.Range("E24").Value = 0.1028055226074426
$ws.Range("C25").Value = 32.97614095206262
$ws.Range("D25").Value = 0.1261409520626202
$ws.Range("E25").Value = 0.01591153978726425
$ws.Range("C26").Value = 33.03574193207213
$ws.Range("D26").Value = 0.1357419320721291
$ws.Range("E26").Value = 0.01842587212267452
$ws.Range("C27").Value = 33.2519920345837
$ws.Range("D27").Value = 0.1519920345836994
$ws.Range("E27").Value = 0.02310157857689247
$ws.Range("C28").Value = 33.5578811770355
$ws.Range("D28").Value = 0.157881177035506
$ws.Range("E28").Value = 0.02492646606211678
$ws.Range("C29").Value = 33.6559370622993
$ws.Range("D29").Value = -0.04406293770070135
$ws.Range("E29").Value = 0.001941542478815889
$ws.Range("C30").Value = 33.75523765042924
$ws.Range("D30").Value = -0.3447623495707646
$ws.Range("E30").Value = 0.1188610776815541
$ws.Range("C31").Value = 34.32940713013046
$ws.Range("D31").Value = -0.07059286986953595
$ws.Range("E31").Value = 0.004983353276417236
$ws.Range("C32").Value = 34.63000944736938
$ws.Range("D32").Value = -0.2699905526306168
$ws.Range("E32").Value = 0.07289489850978585
$ws.Range("C33").Value = 35.70998493816687
$ws.Range("D33").Value = 0.4099849381668719
$ws.Range("E33").Value = 0.1680876495236938
$ws.Range("C34").Value = 35.96754034663636
$ws.Range("D34").Value = 0.2675403466363591
$ws.Range("E34").Value = 0.07157783707830316
$ws.Range("C35").Value = 36.3516462931155
$ws.Range("D35").Value = 0.05164629311550328
$ws.Range("E35").Value = 0.002667339592572482
$ws.Range("C36").Value = 36.87240230435085
$ws.Range("D36").Value = 0.07240230435085238
$ws.Range("E36").Value = 0.005242093675313457
$ws.Range("C37").Value = 37.02096205104752
$ws.Range("D37").Value = -0.2790379489524781
$ws.Range("E37").Value = 0.07786217695560579
$ws.Range("C38").Value = 37.86916685439945
$ws.Range("D38").Value = -0.03083314560055328
$ws.Range("E38").Value = 0.0009506828676249181
$ws.Range("C39").Value = 38.57339737651738
$ws.Range("D39").Value = 0.07339737651738432
$ws.Range("E39").Value = 0.005387174879634678
$ws.Range("C40").Value = 39.13221582775235
$ws.Range("D40").Value = 0.2322158277523556
$ws.Range("E40").Value = 0.05392419065871169
$ws.Range("C41").Value = 39.44683827096657
$ws.Range("D41").Value = 0.04683827096657467
$ws.Range("E41").Value = 0.002193823627138272
$ws.Range("C42").Value = 39.72154326227819
$ws.Range("D42").Value = -0.1784567377218096
$ws.Range("E42").Value = 0.03184680723831074
$ws.Range("C43").Value = 39.76597028438642
$ws.Range("D43").Value = -0.3340297156135819
$ws.Range("E43").Value = 0.1115758509128904
$ws.Range("C44").Value = 39.89940613788547
$ws.Range("D44").Value = -0.7005938621145305
$ws.Range("E44").Value = 0.4908317596325537
$ws.Range("C45").Value = 40.18810999186909
$ws.Range("D45").Value = -0.7118900081309079
$ws.Range("E45").Value = 0.506787383676624
$ws.Range("C46").Value = 41.19510746503966
$ws.Range("D46").Value = -0.004892534960340811
$ws.Range("E46").Value = 0.00002393689833815706
$ws.Range("C47").Value = 41.27281613400426
$ws.Range("D47").Value = -0.227183865995741
$ws.Range("E47").Value = 0.05161250896877081
$ws.Range("C48").Value = 41.706695964021
$ws.Range("D48").Value = -0.09330403597899561
$ws.Range("E48").Value = 0.008705643129969708
$ws.Range("C49").Value = 42.37774458888481
$ws.Range("D49").Value = 0.1777445888848064
$ws.Range("E49").Value = 0.03159313887782884
$ws.Range("C50").Value = 43.50534887802869
$ws.Range("D50").Value = 0.8053488780286884
$ws.Range("E50").Value = 0.6485868153420672
$ws.Range("C51").Value = 43.92497559223219
$ws.Range("D51").Value = 0.224975592232191
$ws.Range("E51").Value = 0.05061401710022508

$ws.Range("C52").Value = 0.1449619343228541
$ws.Range("E52").Value = 4.216651199693445
$ws.Range("E53").Value = 0.08433302399386891

